$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("measure_config")

# Remove the rare_event_chart column (column G); everything to the right shifts left
$ws.Columns.Item(7).Delete()

# New trailing columns for the rare-events flag follow-up fields
$ws.Range("N1").Value = "allowable_days_lag"
$ws.Range("O1").Value = "reviewed_at"
$ws.Range("P1").Value = "escalated_to"

# Give the three new columns sensible custom widths (matches target layout)
$ws.Columns.Item(14).ColumnWidth = 20.666666666666668
$ws.Columns.Item(15).ColumnWidth = 26
$ws.Columns.Item(16).ColumnWidth = 30.166666666666668

# Row 2 (ref 1 - Attendances): rebase demonstration data
$ws.Range("L2").Value = "`"2020-04-27`""
$ws.Range("M2").Value = "Rebased to demonstrate the method.  Add the rebase_dates and rebase_comment to 'measure_config.xlsx'."

# Row 3 (ref 5 - Capacity): clear the old rebase date, add escalation meeting info
$ws.Range("L3").ClearContents()
$ws.Range("O3").Value = "Service performance meeting"
$ws.Range("P3").Value = "Divisional perfomance meeting"

# Update the active cell selection
$ws.Range("A4").Select()
